$p = $ppt.ActivePresentation

# Slide 1: title "Header" + " " + "with" + " " + "inline code" (Consolas)
# -> merge the first four plain runs into a single run "Header with ",
#    keep the Consolas "inline code" run untouched.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$c1 = $tr1.Characters(1, 12)
$c1.Text = "Header with "

# Slide 2: title "Syntax" + " " + "highlighting" -> merge into one run.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$c2 = $tr2.Characters(1, $tr2.Length)
$c2.Text = "Syntax highlighting"

# Slide 3: title "Two" + " " + "column" + " " + "slide" -> merge into one run.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$c3 = $tr3.Characters(1, $tr3.Length)
$c3.Text = "Two column slide"
